$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "2025-04-28 04:42:23"
$ws.Cells.Item(56, 3).Value = "James Davis received Telsla Battery 4 from Suppliers Battery New.`nTelsla Battery 4's state was New.`nThus James Davis carried out the following actions:`nStore, .`nNow James Davis is Tired, feeling that the task was Challenging.`n"
$ws.Cells.Item(56, 3).WrapText = $true

$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = "2025-04-28 04:43:56"
$ws.Cells.Item(57, 3).Value = "James Davis found New Battery from Ford.`nNow James Davis is Excited, feeling that the task was Challenging.`n"
$ws.Cells.Item(57, 3).WrapText = $true

$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = "2025-04-28 04:54:31"
$ws.Cells.Item(58, 3).Value = "James Davis added Brand new battery to the database.`nSerial Number is 58t3952310422.`nPart Number is 34.`nItem Type is 1231.`nLocation is floor space 3.`nNow James Davis is Excited, feeling that the task was Stressful.`n"
$ws.Cells.Item(58, 3).WrapText = $true

$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = "2025-04-28 04:57:23"
$ws.Cells.Item(59, 3).Value = "James Davis moved New Battery from Ford from floor space 2 to floor space 2.`nNow James Davis is Happy, feeling that the task was Fun.`n"
$ws.Cells.Item(59, 3).WrapText = $true

$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = "2025-04-28 05:04:55"
$ws.Cells.Item(60, 3).Value = "James Davis received New Battery from Ford from Suppliers Never Death Row.`nNew Battery from Ford's state was New.`nThus James Davis carried out the following actions:`nUpdate Battery Status, Store, .`nNow James Davis is Excited, feeling that the task was Fun.`n"
$ws.Cells.Item(60, 3).WrapText = $true

$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "2025-04-28 05:11:02"
$ws.Cells.Item(61, 3).Value = "James Davis took picture of Telsla Battery 4.`nNow James Davis is Excited, feeling that the task was Rewarding.`n"
$ws.Cells.Item(61, 3).WrapText = $true

$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "2025-04-28 05:25:32"
$ws.Cells.Item(62, 3).Value = "James Davis added Tesla to the database.`nSerial Number: 87756453234567553`nPart Number: 54`nItem Type: 3`nLocation: shelf space 1`nNow James Davis is Excited, feeling the task was Fun.`n"
$ws.Cells.Item(62, 3).WrapText = $true
